$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: new absence record (Nom / Prenom / Date fin de l'absence)
$ws.Range("A4").Value = "Camur "
$ws.Range("B4").Value = "Abdullah "
$ws.Range("C3").Copy()
$ws.Range("C4").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("C4").Value = (Get-Date -Year 2023 -Month 6 -Day 3 -Hour 0 -Minute 0 -Second 0)

# Clear the leftover single-column names that used to occupy B5:B7
$ws.Range("B5").Value = $null
$ws.Range("B6").Value = $null
$ws.Range("B7").Value = $null

# Zoom in on the active sheet view
$ws.Application.ActiveWindow.Zoom = 225
